$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet lists pharmacy products in rows 4..74 (row N => item N-3), with
# row 75 a subtotal row and row 76 a footer row. We are inserting one new
# product row ("سرنجات انسولين") right before the current row 68
# ("فازلين بيور كبير"), which pushes that item and everything below it
# (through the old subtotal/footer rows) down by one row.
#
# Column A (sequence number) is already simply "row - 3" for every data row,
# so it needs no edits at all for the rows that merely shift down, and we
# only need to set it explicitly for the brand new last data row.
# Column C/D/E/F/G/I/J/K/M are always blank (part of merged B:G / H:K / L:M).
# ---------------------------------------------------------------------------

# 1) Snapshot the current (pre-shift) B/H/L/N values for rows 68-74, top to
#    bottom, before anything gets overwritten.
$rows = 68..74
$snapB = @{}
$snapH = @{}
$snapL = @{}
$snapN = @{}
foreach ($r in $rows) {
    $snapB[$r] = $ws.Range("B$r").Value()
    $snapH[$r] = $ws.Range("H$r").Value()
    $snapL[$r] = $ws.Range("L$r").Value()
    $snapN[$r] = $ws.Range("N$r").Value()
}

# Snapshot the old subtotal (row 75) and footer (row 76) content too.
$oldK75 = $ws.Range("K75").Value()
$oldA76 = $ws.Range("A76").Value()
$oldF76 = $ws.Range("F76").Value()
$oldI76 = $ws.Range("I76").Value()

# 2) Give the brand-new last data row (75) the same per-column formatting
#    already used by every other data row (row 68 has exactly that pattern),
#    before we touch row 68's own content.
$ws.Range("A68:N68").Copy()
$ws.Range("A75:N75").PasteSpecial(-4122)

# 3) Shift the product rows down by one (write bottom-up so sources aren't
#    clobbered before they are read - we already snapshotted them anyway).
foreach ($r in $rows) {
    $dest = $r + 1
    $ws.Range("B$dest").Value = $snapB[$r]
    $ws.Range("H$dest").Value = $snapH[$r]
    $ws.Range("L$dest").Value = $snapL[$r]
    $ws.Range("N$dest").Value = $snapN[$r]
}
$ws.Range("A75").Value = 72

# 4) Write the new product into the now-vacated row 68.
$ws.Range("B68").Value = "سرنجات انسولين"
$ws.Range("H68").Value = "177:0"
$ws.Range("L68").Value = 14
$ws.Range("N68").Value = "2:0"

# 5) Move the subtotal row down to 76 and update the sum to include the
#    new row's price (4745.84 + 14 = 4759.84).
$ws.Range("K76").Value = $oldK75 + 14

# 6) Move the footer row down to 77 (timestamp / page / credit line).
$ws.Range("A77").Value = $oldA76
$ws.Range("F77").Value = $oldF76
$ws.Range("I77").Value = $oldI76

# 7) Fix up row heights: rows 68-74 keep their original heights untouched.
#    The new last data row (75) gets the standard data-row height, the
#    subtotal row (now 76) gets a slightly shorter auto-fit height, and the
#    footer row (now 77) keeps its previous height.
$ws.Rows("75:75").RowHeight = 25.5
$ws.Rows("76:76").RowHeight = 25.5
$ws.Rows("77:77").RowHeight = 16.5

# 8) Rebuild the merged-cell layout for the affected rows: unmerge the old
#    subtotal/footer merges, then re-merge at their new row numbers plus the
#    three new merges needed for row 75's product columns.
$ws.Range("K75:N75").UnMerge()
$ws.Range("A76:E76").UnMerge()
$ws.Range("F76:G76").UnMerge()
$ws.Range("I76:N76").UnMerge()

$ws.Range("B75:G75").Merge()
$ws.Range("H75:K75").Merge()
$ws.Range("L75:M75").Merge()
$ws.Range("K76:N76").Merge()
$ws.Range("A77:E77").Merge()
$ws.Range("F77:G77").Merge()
$ws.Range("I77:N77").Merge()
